$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 27
$ws.Range("B7").Value = "Update index.py"
$ws.Range("C7").Value = "riya-morankar"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "edit1 to main"

# "2025-06-17" looks like a date, so Excel would normally auto-convert it to
# a date serial number. Force it to stay plain text (matching the other
# rows' Date column, which are stored as literal text) by pre-formatting the
# cell as Text, then resetting the style back to Normal afterwards so no
# lingering number-format is left applied to the cell.
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2025-06-17"
$ws.Range("F7").Style = "Normal"
